$d = $word.ActiveDocument

# The document carries a "_GoBack" bookmark (Word's automatic "last edit
# location" marker). It currently sits right after the "Delete an existing
# product from the database." run. This edit moves it into the middle of
# the earlier paragraph, right after "...Lecture Note Code Examples zip
# file " and before "in this week's Learning Resources...", which splits
# that run into two runs around the (empty) bookmark.

# Remove the existing _GoBack bookmark from its old location.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Locate the new split point: right after "...zip file " and before
# "in this week's Learning Resources...".
$searchRange = $d.Content
$found = $searchRange.Find.Execute("Lecture Note Code Examples zip file ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Collapse to an insertion point right after the matched text and add the
# bookmark there.
$insertPoint = $d.Range($searchRange.End, $searchRange.End)
$d.Bookmarks.Add("_GoBack", $insertPoint)
